$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6648
$ws.Range("J3").Value = 7024
$ws.Range("J4").Value = 1526
$ws.Range("J6").Value = 9335
$ws.Range("J7").Value = 25077

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 431
$ws.Range("J7").Value = 1573

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 147
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 501

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 264
$ws.Range("J6").Value = 396
$ws.Range("J7").Value = 1134

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 226
$ws.Range("J3").Value = 261
$ws.Range("J4").Value = 28
$ws.Range("J6").Value = 226
$ws.Range("J7").Value = 771

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 200
$ws.Range("J7").Value = 729
$ws.Range("J8").Value = 1573
$ws.Range("J11").Value = 429
$ws.Range("J12").Value = 53
$ws.Range("J18").Value = 212
$ws.Range("J19").Value = 733
$ws.Range("J20").Value = 523
$ws.Range("J23").Value = 229
$ws.Range("J24").Value = 77
$ws.Range("J25").Value = 124
$ws.Range("J27").Value = 149
$ws.Range("J29").Value = 1370
$ws.Range("J31").Value = 246
$ws.Range("J32").Value = 41
$ws.Range("J33").Value = 1134
$ws.Range("J36").Value = 338
$ws.Range("J37").Value = 771
$ws.Range("J42").Value = 1080
$ws.Range("J44").Value = 191
$ws.Range("J46").Value = 84
$ws.Range("J48").Value = 283
$ws.Range("J49").Value = 160
$ws.Range("J51").Value = 307
$ws.Range("J52").Value = 635
$ws.Range("J54").Value = 479
$ws.Range("J61").Value = 28
$ws.Range("I63").Value = 180
$ws.Range("J64").Value = 167
$ws.Range("J66").Value = 76
$ws.Range("J67").Value = 945
$ws.Range("J69").Value = 55
$ws.Range("I79").Value = 748
$ws.Range("J79").Value = 707
$ws.Range("J83").Value = 501
$ws.Range("J85").Value = 1043
$ws.Range("J86").Value = 160
$ws.Range("J87").Value = 84
$ws.Range("J88").Value = 259
$ws.Range("J89").Value = 320
$ws.Range("J90").Value = 267
$ws.Range("J91").Value = 289
$ws.Range("J94").Value = 267
$ws.Range("J96").Value = 275
$ws.Range("J101").Value = 25077

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 246

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 242
$ws.Range("J3").Value = 351
$ws.Range("J7").Value = 945

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 225
$ws.Range("J7").Value = 479

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 483
$ws.Range("J6").Value = 346
$ws.Range("J7").Value = 1370

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 139
$ws.Range("J7").Value = 283

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J4").Value = 36
$ws.Range("J7").Value = 733

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 61
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 231
$ws.Range("J6").Value = 574
$ws.Range("J7").Value = 1080

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 101

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 78
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 74
$ws.Range("J7").Value = 275

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 240
$ws.Range("I4").Value = 39
$ws.Range("J6").Value = 210
$ws.Range("I7").Value = 748
$ws.Range("J7").Value = 707

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 176
$ws.Range("J7").Value = 523

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 108
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 225
$ws.Range("J3").Value = 222
$ws.Range("J4").Value = 28
$ws.Range("J7").Value = 729

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 52
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J4").Value = 25
$ws.Range("J7").Value = 429

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 200

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 54
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 82
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 280
$ws.Range("J6").Value = 303
$ws.Range("J7").Value = 1043

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 147
$ws.Range("J4").Value = 24
$ws.Range("J7").Value = 635

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 28
